# Revert "RESTORE: Recover all 973 original multi-industry template files"
# i.e. rename the "Artificial Intelligence and Machine Learning" themed
# template text back to "Product Development", and restore the blank
# spacer rows that existed in the original layout.

$wb = $excel.ActiveWorkbook

function Stamp-EmptyRow($ws, $rowNum) {
    # Forces openpyxl/Excel to persist a truly-empty <row r="N"/> element
    # (no cells, no extra attributes) in the saved sheetData.
    $ws.Rows.Item($rowNum).OutlineLevel = 0
}

# ---------------------------------------------------------------------
# Sheet 1: Instructions & User Guide
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Instructions & User Guide")

$ws1.Range("A1").Value = "Product Development Comprehensive Budget - User Guide & Instructions"
$ws1.Range("A56").Value = "📋 PRODUCT DEVELOPMENT PROJECT OVERVIEW"
$ws1.Range("B59").Value = "Data Scientists, Product Engineers, Product Architects, DevOps Engineers..."

foreach ($r in 10, 20, 28, 37, 45, 54, 55, 60) {
    Stamp-EmptyRow $ws1 $r
}

# ---------------------------------------------------------------------
# Sheet 2: Budget Summary
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Budget Summary")

$ws2.Range("A1").Value = "Product Development - Executive Budget Summary"

foreach ($r in 2, 6) {
    Stamp-EmptyRow $ws2 $r
}

# ---------------------------------------------------------------------
# Sheet 3: Resources
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Resources")

$ws3.Range("A1").Value = "Product Development - Resources Budget"
$ws3.Range("A5").Value = "Product Engineers"
$ws3.Range("A6").Value = "Product Architects"

foreach ($r in 2, 11) {
    Stamp-EmptyRow $ws3 $r
}

# ---------------------------------------------------------------------
# Sheet 4: Logistics
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Logistics")

$ws4.Range("A1").Value = "Product Development - Logistics Budget"

foreach ($r in 2, 9) {
    Stamp-EmptyRow $ws4 $r
}

# ---------------------------------------------------------------------
# Sheet 5: Technology
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Technology")

$ws5.Range("A1").Value = "Product Development - Technology Budget"
$ws5.Range("A5").Value = "Product Platform Licenses"

foreach ($r in 2, 10) {
    Stamp-EmptyRow $ws5 $r
}

# ---------------------------------------------------------------------
# Sheet 6: Training
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Training")

$ws6.Range("A1").Value = "Product Development - Training Budget"
$ws6.Range("A4").Value = "Product Development Certification Programs"
$ws6.Range("A10").Value = "TOTAL TRProductNING"

foreach ($r in 2, 9) {
    Stamp-EmptyRow $ws6 $r
}

# ---------------------------------------------------------------------
# Sheet 7: Contingency
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Contingency")

$ws7.Range("A1").Value = "Product Development - Contingency Budget"

foreach ($r in 2, 5, 11, 13) {
    Stamp-EmptyRow $ws7 $r
}

# ---------------------------------------------------------------------
# Sheet 8: Timeline
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("Timeline")

$ws8.Range("A1").Value = "Product Development - Budget Timeline"

foreach ($r in 2) {
    Stamp-EmptyRow $ws8 $r
}
